# ---------------------------------------------------------------------------
# maldiims-metadata.xlsx: add a "version"/"description" pair of columns at the
# front of the "Export as TSV" sheet, backed by a new "version list" sheet.
# ---------------------------------------------------------------------------

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Export as TSV")

# ---------------------------------------------------------------------------
# 1. Collect the existing header-row comments (A1..AD1) BEFORE we touch
#    anything, keyed by their current column letter, so we can re-create them
#    two columns to the right after the insert (comments do not travel with
#    an EntireColumn.Insert()).
# ---------------------------------------------------------------------------
$oldCols = @("A","B","C","D","E","F","G","H","I","J","K","L","M","N","O","P","Q","R","S","T","U","V","W","X","Y","Z","AA","AB","AC","AD")
$oldComments = @{}
foreach ($col in $oldCols) {
    $cell = $ws.Range($col + "1")
    if ($cell.Comment -ne $null) {
        $oldComments[$col] = $cell.Comment.Text()
    }
}

# ---------------------------------------------------------------------------
# 2. Insert two new blank columns at the front (A:B). Everything that used to
#    live in A..AD now lives in C..AF; Excel automatically re-targets the
#    data-validation sqrefs and formulas, which is exactly what the target
#    workbook shows.
# ---------------------------------------------------------------------------
$ws.Range("A:B").EntireColumn.Insert()

# ---------------------------------------------------------------------------
# 3. Helper: convert a 1-based column index to its letter(s).
# ---------------------------------------------------------------------------
function Get-ColLetter($idx) {
    $s = ""
    while ($idx -gt 0) {
        $rem = ($idx - 1) % 26
        $s = [char](65 + $rem) + $s
        $idx = [int](($idx - $rem - 1) / 26)
    }
    return $s
}

function Get-ColIndex($letters) {
    $idx = 0
    foreach ($ch in $letters.ToCharArray()) {
        $idx = $idx * 26 + ([int][char]$ch - 64)
    }
    return $idx
}

# ---------------------------------------------------------------------------
# 4. Re-create every old comment two columns over (its new home), then add
#    the two brand-new ones for the inserted "version" / "description"
#    columns.
# ---------------------------------------------------------------------------
foreach ($col in $oldCols) {
    if ($oldComments.ContainsKey($col)) {
        $newCol = Get-ColLetter (Get-ColIndex($col) + 2)
        $cell = $ws.Range($newCol + "1")
        if ($cell.Comment -ne $null) {
            $cell.Comment.Delete()
        }
        $cell.AddComment($oldComments[$col])
    }
}

$ws.Range("A1").AddComment("Version of the schema to use when validating this metadata.")
$ws.Range("B1").AddComment("Free-text description of this assay.")

# ---------------------------------------------------------------------------
# 5. Fill in the new header cells.
# ---------------------------------------------------------------------------
$ws.Range("A1").Value = "version"
$ws.Range("B1").Value = "description"

# ---------------------------------------------------------------------------
# 6. Add a new worksheet, "version list", right after "Export as TSV", and
#    give it its single value (stored as text, like the other *_list sheets).
# ---------------------------------------------------------------------------
$verList = $wb.Worksheets.Add($null, $ws)
$verList.Name = "version list"
$verList.Range("A1").NumberFormat = "@"
$verList.Range("A1").Value = "1"

# ---------------------------------------------------------------------------
# 7. Add the data validation for the new "version" column, matching the style
#    of the other list-backed validations on this sheet.
# ---------------------------------------------------------------------------
$verRange = $ws.Range("A2:A1048576")
$verRange.Validation.Add(3, 1, 1, "='version list'!`$A`$1:`$A`$1")
$verRange.Validation.ErrorTitle = "Value must come from list"
$verRange.Validation.ErrorMessage = "Value must be one of: 1."
$verRange.Validation.ShowInput = $true
$verRange.Validation.ShowError = $true
$verRange.Validation.IgnoreBlank = $true
